# Historico.xlsx — add a new "Homicidios" column before the existing
# "ratio_turistas_residentes" column (R), shifting that column to S,
# then populate the new column with the Homicidios figures that are
# available (years 1995-1997, 2003-2004, 2010-2013, 2015, 2017-2020).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at R; this shifts the former column R
# ("ratio_turistas_residentes", including its header and all data) one
# place to the right, to column S.
$ws.Range("R1").EntireColumn.Insert()

# Header for the freshly inserted column.
$ws.Range("R1").Value = "Homicidios"

# New Homicidios data points.
$ws.Range("R37").Value = 6.94455940999813
$ws.Range("R38").Value = 7.37742376484379
$ws.Range("R39").Value = 7.45770271278924
$ws.Range("R45").Value = 7.48997932167156
$ws.Range("R46").Value = 7.6778416479797
$ws.Range("R52").Value = 8.68662061022504
$ws.Range("R53").Value = 9.013810683079781
$ws.Range("R54").Value = 8.637601847104291
$ws.Range("R55").Value = 7.57207576168001
$ws.Range("R57").Value = 7.20094334280093
$ws.Range("R58").Value = 6.38934150079073
$ws.Range("R59").Value = 5.65469286500665
$ws.Range("R60").Value = 4.75998593558216
$ws.Range("R61").Value = 4.59485426006108
$ws.Range("R62").Value = 3.74921291744897
